$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 83918.914
$ws.Range("I2").Value = 363.83334
$ws.Range("J2").Value = 167474
$ws.Range("K2").Value = 363.83334
$ws.Range("L2").Value = 167474
$ws.Range("M2").Value = -250.83334
$ws.Range("N2").Value = -167700
$ws.Range("H15").Value = 1151.7273
$ws.Range("I15").Value = 1151.7273
$ws.Range("K15").Value = 3455.1819
$ws.Range("M15").Value = -3286.1819
$ws.Range("H29").Value = 3059.6
$ws.Range("I29").Value = 1466.3334
$ws.Range("J29").Value = 5449.5
$ws.Range("K29").Value = 4399.0002
$ws.Range("L29").Value = 16348.5
$ws.Range("M29").Value = -4118.0002
$ws.Range("N29").Value = -16910.5
$ws.Range("H38").Value = 8179.8335
$ws.Range("I38").Value = 27.166666
$ws.Range("J38").Value = 16332.5
$ws.Range("K38").Value = 81.49999800000001
$ws.Range("L38").Value = 48997.5
$ws.Range("M38").Value = 290.500002
$ws.Range("N38").Value = -49741.5
$ws.Range("H99").Value = 311
$ws.Range("J99").Value = 500
$ws.Range("L99").Value = 1500
$ws.Range("N99").Value = -4496
$ws.Range("H100").Value = 3139.3635
$ws.Range("I100").Value = 3153.3
$ws.Range("K100").Value = 3153.3
$ws.Range("M100").Value = -2612.3
$ws.Range("H116").Value = 7301.5
$ws.Range("I116").Value = 7221.75
$ws.Range("K116").Value = 7221.75
$ws.Range("M116").Value = -3779.75
$ws.Range("H125").Value = 1259.7084
$ws.Range("I125").Value = 1358.6666
$ws.Range("J125").Value = 1226.7222
$ws.Range("K125").Value = 12227.9994
$ws.Range("L125").Value = 11040.4998
$ws.Range("M125").Value = -9767.999400000001
$ws.Range("N125").Value = -15960.4998
$ws.Range("H126").Value = 54000
$ws.Range("J126").Value = 54000
$ws.Range("L126").Value = 54000
$ws.Range("N126").Value = -63880

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2918
$ws.Range("I2").Value = 2918
$ws.Range("K2").Value = 2918
$ws.Range("M2").Value = -2805
$ws.Range("H32").Value = 17522.596
$ws.Range("I32").Value = 14961.613
$ws.Range("J32").Value = 57217.832
$ws.Range("K32").Value = 14961.613
$ws.Range("L32").Value = 57217.832
$ws.Range("M32").Value = -14674.613
$ws.Range("N32").Value = -57791.832
$ws.Range("H45").Value = 1245.6774
$ws.Range("I45").Value = 1179.0714
$ws.Range("K45").Value = 1179.0714
$ws.Range("M45").Value = -802.0714
$ws.Range("H63").Value = 4348.294
$ws.Range("I63").Value = 3070.8462
$ws.Range("J63").Value = 8500
$ws.Range("K63").Value = 3070.8462
$ws.Range("L63").Value = 8500
$ws.Range("M63").Value = -2384.8462
$ws.Range("N63").Value = -9872
$ws.Range("H66").Value = 4348.294
$ws.Range("I66").Value = 3070.8462
$ws.Range("J66").Value = 8500
$ws.Range("K66").Value = 15354.231
$ws.Range("L66").Value = 42500
$ws.Range("M66").Value = -11922.231
$ws.Range("N66").Value = -49364
$ws.Range("H88").Value = 1843.4445
$ws.Range("I88").Value = 1361.6666
$ws.Range("J88").Value = 2084.3333
$ws.Range("K88").Value = 1361.6666
$ws.Range("L88").Value = 2084.3333
$ws.Range("M88").Value = -955.6666
$ws.Range("N88").Value = -2896.3333
$ws.Range("H91").Value = 1843.4445
$ws.Range("I91").Value = 1361.6666
$ws.Range("J91").Value = 2084.3333
$ws.Range("K91").Value = 1361.6666
$ws.Range("L91").Value = 2084.3333
$ws.Range("M91").Value = 42.33339999999998
$ws.Range("N91").Value = -4892.3333
$ws.Range("H102").Value = 385995.47
$ws.Range("I102").Value = 455861.88
$ws.Range("K102").Value = 455861.88
$ws.Range("M102").Value = -454239.88
$ws.Range("H116").Value = 2918
$ws.Range("I116").Value = 2918
$ws.Range("K116").Value = 2918
$ws.Range("M116").Value = -624
$ws.Range("H119").Value = 60690.5
$ws.Range("J119").Value = 60690.5
$ws.Range("L119").Value = 60690.5
$ws.Range("N119").Value = -70366.5
$ws.Range("H122").Value = 4155.625
$ws.Range("I122").Value = 3508.5557
$ws.Range("K122").Value = 10525.6671
$ws.Range("M122").Value = -8075.667099999999

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2918
$ws.Range("I3").Value = 2918
$ws.Range("K3").Value = 2918
$ws.Range("M3").Value = -2804
$ws.Range("H12").Value = 6500
$ws.Range("I12").Value = 3000
$ws.Range("K12").Value = 3000
$ws.Range("M12").Value = -2832
$ws.Range("H94").Value = 1131.3846
$ws.Range("I94").Value = 876.8823
$ws.Range("K94").Value = 876.8823
$ws.Range("M94").Value = -425.8823
$ws.Range("H99").Value = 2765.75
$ws.Range("I99").Value = 2311.25
$ws.Range("K99").Value = 2311.25
$ws.Range("M99").Value = -813.25

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 874
$ws.Range("I22").Value = 874
$ws.Range("K22").Value = 874
$ws.Range("M22").Value = -524
$ws.Range("H31").Value = 3487.9856
$ws.Range("I31").Value = 3103.9443
$ws.Range("J31").Value = 3623.5293
$ws.Range("K31").Value = 3103.9443
$ws.Range("L31").Value = 3623.5293
$ws.Range("M31").Value = -2808.9443
$ws.Range("N31").Value = -4213.5293
$ws.Range("H34").Value = 3487.9856
$ws.Range("I34").Value = 3103.9443
$ws.Range("J34").Value = 3623.5293
$ws.Range("K34").Value = 3103.9443
$ws.Range("L34").Value = 3623.5293
$ws.Range("M34").Value = -2901.9443
$ws.Range("N34").Value = -4027.5293
$ws.Range("H62").Value = 10648.728
$ws.Range("J62").Value = 10780.5
$ws.Range("L62").Value = 10780.5
$ws.Range("N62").Value = -12028.5
$ws.Range("H65").Value = 10648.728
$ws.Range("J65").Value = 10780.5
$ws.Range("L65").Value = 53902.5
$ws.Range("N65").Value = -60142.5
$ws.Range("H105").Value = 1113.174
$ws.Range("I105").Value = 1085.4615
$ws.Range("K105").Value = 1085.4615
$ws.Range("M105").Value = 661.5385000000001
$ws.Range("H110").Value = 135500
$ws.Range("J110").Value = 135500
$ws.Range("L110").Value = 135500
$ws.Range("N110").Value = -143680

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 6456.625
$ws.Range("J3").Value = 5775.5
$ws.Range("L3").Value = 5775.5
$ws.Range("N3").Value = -6007.5
$ws.Range("H4").Value = 8291.5
$ws.Range("I4").Value = 5000
$ws.Range("J4").Value = 8949.799999999999
$ws.Range("K4").Value = 5000
$ws.Range("L4").Value = 8949.799999999999
$ws.Range("M4").Value = -4888
$ws.Range("N4").Value = -9173.799999999999
$ws.Range("H80").Value = 7441.0347
$ws.Range("I80").Value = 2451.625
$ws.Range("J80").Value = 9341.762000000001
$ws.Range("K80").Value = 2451.625
$ws.Range("L80").Value = 9341.762000000001
$ws.Range("M80").Value = -1453.625
$ws.Range("N80").Value = -11337.762
$ws.Range("H83").Value = 7441.0347
$ws.Range("I83").Value = 2451.625
$ws.Range("J83").Value = 9341.762000000001
$ws.Range("K83").Value = 12258.125
$ws.Range("L83").Value = 46708.81
$ws.Range("M83").Value = -7266.125
$ws.Range("N83").Value = -56692.81
$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("M93").ClearContents()
$ws.Range("H97").Value = 1258.1428
$ws.Range("I97").Value = 599.5
$ws.Range("K97").Value = 599.5
$ws.Range("M97").Value = -103.5
$ws.Range("H102").Value = 20342.65
$ws.Range("I102").Value = 23900.746
$ws.Range("J102").Value = 9134.65
$ws.Range("K102").Value = 23900.746
$ws.Range("L102").Value = 9134.65
$ws.Range("M102").Value = -22278.746
$ws.Range("N102").Value = -12378.65
$ws.Range("H105").Value = 55990.184
$ws.Range("J105").Value = 55990.184
$ws.Range("L105").Value = 55990.184
$ws.Range("N105").Value = -62978.184

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 50746.75
$ws.Range("I40").Value = 58153.258
$ws.Range("J40").Value = 21943.666
$ws.Range("K40").Value = 58153.258
$ws.Range("L40").Value = 21943.666
$ws.Range("M40").Value = -58017.258
$ws.Range("N40").Value = -22215.666
$ws.Range("H122").Value = 4348.7427
$ws.Range("I122").Value = 4001.2693
$ws.Range("K122").Value = 12003.8079
$ws.Range("M122").Value = -9553.8079

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 9991.6
$ws.Range("I122").Value = 9990.777
$ws.Range("K122").Value = 29972.331
$ws.Range("M122").Value = -27522.331
$ws.Range("H126").Value = 34590.848
$ws.Range("I126").Value = 43448.1
$ws.Range("K126").Value = 130344.3
$ws.Range("M126").Value = -127874.3
